# Re-pull / push updated data for the dSF column (F) on Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 5
$ws.Range("F3").Value = -9
$ws.Range("F6").Value = 5
$ws.Range("F12").Value = 7
$ws.Range("F13").Value = 0
